$wb = $excel.ActiveWorkbook
$eda = $wb.Worksheets.Item("EDA")
$new = $wb.Worksheets.Add($null, $eda)
$new.Name = "Alt1"

$cells = @("C1","D1","E1","G1","H1","I1","N1","C2","D2","E2","G2","H2","I2","C3","D3","E3","G3","H3","I3","O3","C4","D4","E4","G4","H4","I4","O4","C5","D5","E5","G5","H5","I5","O5","C6","D6","E6","G6","H6","I6","C7","D7","E7","G7","C8","D8","E8","G8","C9","D9","E9","G9","C10","D10","E10","G10","C11","D11","E11","G11","C12","D12","E12","C13","D13","E13","C14","D14","E14","C15","D15","E15","C16","D16","E16","C17","D17","E17","C18","D18","E18","C19","D19","E19","C20","D20","E20","C21","D21","E21","C22","D22","E22","C23","D23","E23","C24","D24","E24","C25","D25","E25","C26","D26","E26","C27","D27","E27")
foreach ($addr in $cells) {
    $s = $eda.Range($addr)
    $s.Copy()
    $d = $new.Range($addr)
    $d.PasteSpecial(-4122)
    $d.Value2 = $s.Value2
}
Write-Host "UsedRange:" $new.UsedRange.Address()
$new.Range("C1:E1").Merge()
$new.Range("G1:I1").Merge()
